$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add data for columns C, D, E for rows 1-5
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 3

$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 4

$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 5

$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 2

$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 1

# Update selection to D6
$ws.Range("D6").Select()
